# Auto-update draw results: append the 2025-10-07 Pick 3 draw as a new
# row (21) at the bottom of the results table, then grow the table's
# dimension accordingly.
#
# All existing cells in the sheet are plain text (numeric-looking values
# such as dates and the "251007" phase code are intentionally stored as
# text, which is why the sheet carries a numberStoredAsText ignoredError
# exemption). To keep the appended row consistent with that convention we
# force the new cells to Text format before writing the values so Excel
# does not auto-convert the date/number-looking strings into a real date
# serial / number, then clear the formatting again so the new row doesn't
# pick up a distinct cell style from the rest of the (unstyled) table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 21

$ws.Range("A" + $newRow + ":E" + $newRow).NumberFormat = "@"

$ws.Range("A" + $newRow).Value = "2025-10-07"
$ws.Range("B" + $newRow).Value = "Pick 3"
$ws.Range("C" + $newRow).Value = "251007"
$ws.Range("D" + $newRow).Value = "1-6-3"
$ws.Range("E" + $newRow).Value = "2025-10-07T21:37:49.707+04:00"

$ws.Range("A" + $newRow + ":E" + $newRow).ClearFormats()

# Grow the "numbers stored as text" error-check suppression to cover the
# freshly appended row too (mirrors the sheet's existing A1:E20 -> A1:E21
# range growth for that ignored-error entry).
$ws.Range("A1:E" + $newRow).Errors.Item(3).Ignore = $true
